$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cellRef, $val) {
    $rng = $ws.Range($cellRef)
    $rng.NumberFormat = "@"
    $rng.Value = $val
    $rng.ClearFormats()
}

Set-TextValue "D2" "25.939.73"
$ws.Range("E2").Value = "  -0.40%  "

Set-TextValue "D3" "1.621.74"
$ws.Range("E3").Value = "  -0.96%  "

$ws.Range("E4").Value = "  -0.25%  "

Set-TextValue "D5" "212.58"
$ws.Range("E5").Value = "  -0.98%  "

$ws.Range("B6").Value = "USDC"
$ws.Range("C6").Value = "https://coinranking.com/coin/aKzUVe4Hh_CON+usdc-usdc"
Set-TextValue "D6" "1.00"
$ws.Range("E6").Value = "  -0.25%  "

$ws.Range("B7").Value = "XRP"
$ws.Range("C7").Value = "https://coinranking.com/coin/-l8Mn2pVlRs-p+xrp-xrp"
Set-TextValue "D7" "0.495"
$ws.Range("E7").Value = "  -1.94%  "

$ws.Range("E8").Value = "  -0.91%  "

Set-TextValue "D9" "0.0621"
$ws.Range("E9").Value = "  -1.03%  "

$ws.Range("E10").Value = "  -1.70%  "

Set-TextValue "D11" "0.0792"
$ws.Range("E11").Value = "  -0.36%  "

Set-TextValue "D12" "1.846.95"
$ws.Range("E12").Value = "  -1.02%  "

$ws.Range("B13").Value = "WrappedEther"
$ws.Range("C13").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
Set-TextValue "D13" "1.628.10"
$ws.Range("E13").Value = "  -2.53%  "

$ws.Range("B14").Value = "Polkadot"
$ws.Range("C14").Value = "https://coinranking.com/coin/25W7FG7om+polkadot-dot"
Set-TextValue "D14" "4.15"
$ws.Range("E14").Value = "  -1.58%  "

$ws.Range("E15").Value = "  -1.53%  "

Set-TextValue "D16" "25.957.02"

Set-TextValue "D17" "61.74"
$ws.Range("E17").Value = "  -0.87%  "

Set-TextValue "D18" "0.0₃0740"
$ws.Range("E18").Value = "  -1.03%  "

$ws.Range("E19").Value = "  -0.28%  "

Set-TextValue "D20" "192.19"
$ws.Range("E20").Value = "  +0.45%  "

$ws.Range("E21").Value = "  -0.58%  "

Set-TextValue "D22" "9.55"
$ws.Range("E22").Value = "  -0.80%  "

$ws.Range("E23").Value = "  -2.31%  "

$ws.Range("E24").Value = "  +0.60%  "

Set-TextValue "D25" "144.18"
$ws.Range("E25").Value = "  +0.26%  "

Set-TextValue "D26" "1.01"
$ws.Range("E26").Value = "  -0.29%  "

$ws.Range("E27").Value = "  -3.57%  "

$ws.Range("E28").Value = "  -1.92%  "

Set-TextValue "D29" "15.23"
$ws.Range("E29").Value = "  -0.27%  "

$ws.Range("E30").Value = "  -1.19%  "

Set-TextValue "D31" "0.0479"
$ws.Range("E31").Value = "  -1.43%  "

$ws.Range("E32").Value = "  -1.51%  "

Set-TextValue "D33" "3.10"
$ws.Range("E33").Value = "  -2.58%  "

$ws.Range("E34").Value = "  -0.70%  "

$ws.Range("E35").Value = "  -1.22%  "

Set-TextValue "D36" "1.127.44"
$ws.Range("E36").Value = "  -0.03%  "

Set-TextValue "D37" "0.846"
$ws.Range("E37").Value = "  -3.65%  "

$ws.Range("E38").Value = "  -1.76%  "

Set-TextValue "D39" "0.518"
$ws.Range("E39").Value = "  -1.98%  "

$ws.Range("E40").Value = "  -1.45%  "

Set-TextValue "D41" "97.88"
$ws.Range("E41").Value = "  -1.18%  "

Set-TextValue "D42" "1.758.05"
$ws.Range("E42").Value = "  -0.81%  "

Set-TextValue "D43" "0.757"
$ws.Range("E43").Value = "  -3.84%  "

Set-TextValue "D44" "5.14"
$ws.Range("E44").Value = "  -3.09%  "

$ws.Range("E45").Value = "  -0.74%  "

Set-TextValue "D46" "1.51"
$ws.Range("E46").Value = "  +1.82%  "

Set-TextValue "D47" "54.15"
$ws.Range("E47").Value = "  -2.61%  "

$ws.Range("E48").Value = "  -1.70%  "

$ws.Range("E50").Value = "  -1.66%  "

$ws.Range("E51").Value = "  -0.08%  "
